# "add n/a in blank cells"
# The two blank cells in the last data row (B7/C7) get a literal "n/a"
# value, and the sheet's view/selection moves to H7 (scrolled down a bit)
# to match where the user clicked after typing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "n/a"
$ws.Range("C7").Value = "n/a"

# Reflect the scrolled viewport / new selection left by the edit.
try {
    $excel.ActiveWindow.ScrollRow = 5
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # Best effort only - some hosts don't expose window scroll state.
}

$ws.Range("H7").Select()
